# Add a new spell ("Feign Death") as a new row at the bottom of the spell
# table on Sheet1, matching the "Some new spells added" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

$ws.Cells.Item($row, 1).Value2  = 'Feign Death'
$ws.Cells.Item($row, 2).Value2  = 'Healing'
$ws.Cells.Item($row, 3).Value2  = 'fautis'
$ws.Cells.Item($row, 4).Value2  = 'Ritual (30 minutes)'
$ws.Cells.Item($row, 6).Value2  = '(1+$2\times$PP) hours'
$ws.Cells.Item($row, 7).Value2  = 3
$ws.Cells.Item($row, 10).Value2 = 'When cast upon a willing living being, they are placed into a state of suspended animation which perfectly replicates the outward appearance of death. Divination checks with a casting check greater than this spell may peirce the deception. The target is blinded, deafened and physically incapacitated for the duration of the spell. The caster may revive them as a minor action. '

# Row 40 is tall (wrapped Effect text), like the other multi-line rows above it.
$ws.Rows.Item($row).RowHeight = 37.3

# Column D ("Type") needs to widen to fit "Ritual (30 minutes)" - it used to
# share a width with column E, now it gets its own, wider, one.
$ws.Columns.Item(4).ColumnWidth = 12.75

# Leave the selection on the newly-added Effect cell, one row below the new
# data (matches where the author's cursor ended up after typing the row).
$ws.Range("J41").Select()
